# Scene.xlsx config update:
# "Share" column (K) controls whether players are grouped into the same
# scene instance. For the cloned village scenes this must be turned off
# (0/FALSE) so a new group is created per payer login instead of sharing
# an existing clone's group.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(10, 11, 13, 14, 15)
foreach ($row in $rows) {
    $ws.Cells.Item($row, 11).Value = 0
}

$ws.Range("K10").Select()
